$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $range = $ws.Range($ref)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextCell $ws 'D2' '98.210.31'
Set-TextCell $ws 'E2' '  -0.34%  '

Set-TextCell $ws 'D3' '3.415.37'
Set-TextCell $ws 'E3' '  +1.00%  '

Set-TextCell $ws 'E4' '  +0.05%  '

Set-TextCell $ws 'D5' '255.52'
Set-TextCell $ws 'E5' '  -1.20%  '

Set-TextCell $ws 'D6' '683.40'
Set-TextCell $ws 'E6' '  +2.14%  '

Set-TextCell $ws 'D7' '1.45'
Set-TextCell $ws 'E7' '  -6.70%  '

Set-TextCell $ws 'D8' '0.434'
Set-TextCell $ws 'E8' '  -5.44%  '

Set-TextCell $ws 'E9' '  -3.45%  '

Set-TextCell $ws 'E10' '  +0.02%  '

Set-TextCell $ws 'D11' '3.413.10'
Set-TextCell $ws 'E11' '  +1.03%  '

Set-TextCell $ws 'D12' '0.216'
Set-TextCell $ws 'E12' '  +3.04%  '

Set-TextCell $ws 'D13' '41.94'
Set-TextCell $ws 'E13' '  -1.50%  '

Set-TextCell $ws 'D14' '6.37'
Set-TextCell $ws 'E14' '  +13.29%  '

Set-TextCell $ws 'D15' '97.962.85'
Set-TextCell $ws 'E15' '  +0.70%  '

Set-TextCell $ws 'D16' '0.0000267'
Set-TextCell $ws 'E16' '  -0.12%  '

Set-TextCell $ws 'D17' '4.043.21'
Set-TextCell $ws 'E17' '  +0.89%  '

Set-TextCell $ws 'D18' '9.11'
Set-TextCell $ws 'E18' '  +19.06%  '

Set-TextCell $ws 'D19' '3.421.08'
Set-TextCell $ws 'E19' '  +1.47%  '

Set-TextCell $ws 'D20' '0.581'
Set-TextCell $ws 'E20' '  +29.03%  '

Set-TextCell $ws 'D21' '17.62'
Set-TextCell $ws 'E21' '  +3.57%  '

Set-TextCell $ws 'D22' '11.10'
Set-TextCell $ws 'E22' '  +5.34%  '

Set-TextCell $ws 'E23' '  -3.99%  '

Set-TextCell $ws 'D24' '510.31'
Set-TextCell $ws 'E24' '  -3.86%  '

Set-TextCell $ws 'D25' '0.0000206'
Set-TextCell $ws 'E25' '  -3.80%  '

Set-TextCell $ws 'D26' '6.59'
Set-TextCell $ws 'E26' '  +4.34%  '

Set-TextCell $ws 'D27' '101.01'
Set-TextCell $ws 'E27' '  -1.22%  '

Set-TextCell $ws 'D28' '12.83'
Set-TextCell $ws 'E28' '  +1.12%  '

Set-TextCell $ws 'D29' '3.603.10'
Set-TextCell $ws 'E29' '  +1.13%  '

Set-TextCell $ws 'E30' '  +0.04%  '

Set-TextCell $ws 'D31' '11.68'
Set-TextCell $ws 'E31' '  +4.09%  '

Set-TextCell $ws 'B32' 'Dai'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D32' '0.998'
Set-TextCell $ws 'E32' '  -0.11%  '

Set-TextCell $ws 'B33' 'Cronos'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D33' '0.196'
Set-TextCell $ws 'E33' '  +3.29%  '

Set-TextCell $ws 'D34' '2.66'
Set-TextCell $ws 'E34' '  +24.44%  '

Set-TextCell $ws 'E35' '  +6.40%  '

Set-TextCell $ws 'D36' '0.997'
Set-TextCell $ws 'E36' '  +1.91%  '

Set-TextCell $ws 'D37' '29.93'
Set-TextCell $ws 'E37' '  +0.29%  '

Set-TextCell $ws 'B38' 'Fetch.AI'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D38' '1.55'
Set-TextCell $ws 'E38' '  +15.30%  '

Set-TextCell $ws 'B39' 'RenderToken'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws 'D39' '8.10'
Set-TextCell $ws 'E39' '  +1.86%  '

Set-TextCell $ws 'D40' '534.42'
Set-TextCell $ws 'E40' '  +1.23%  '

Set-TextCell $ws 'E41' '  -4.41%  '

Set-TextCell $ws 'E42' '  +0.00%  '

Set-TextCell $ws 'D43' '0.881'
Set-TextCell $ws 'E43' '  +5.05%  '

Set-TextCell $ws 'E44' '  +0.08%  '

Set-TextCell $ws 'B45' 'Cosmos'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D45' '9.06'
Set-TextCell $ws 'E45' '  +13.77%  '

Set-TextCell $ws 'B46' 'VeChain'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D46' '0.0437'
Set-TextCell $ws 'E46' '  -3.82%  '

Set-TextCell $ws 'D47' '1.76'
Set-TextCell $ws 'E47' '  +15.68%  '

Set-TextCell $ws 'D48' '3.77'
Set-TextCell $ws 'E48' '  -0.62%  '

Set-TextCell $ws 'D49' '5.81'
Set-TextCell $ws 'E49' '  +12.94%  '

Set-TextCell $ws 'D50' '55.99'
Set-TextCell $ws 'E50' '  +10.12%  '

Set-TextCell $ws 'D51' '3.22'
Set-TextCell $ws 'E51' '  -4.22%  '
